# Update cryptocurrency price/volume data per the 2023-02-15 symbol-list refresh.
# Values are stored as literal text (not numbers/percentages), so we force each
# touched cell to Text format before assigning, then restore the default 'Normal'
# style afterward so no stray number formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "297.83"
Set-TextValue "E2" "2.07%"
Set-TextValue "D3" "42.08"
Set-TextValue "E3" "4.61%"
Set-TextValue "D4" "5.008"
Set-TextValue "E4" "0.04%"
Set-TextValue "D5" "0.07523"
Set-TextValue "E5" "3.19%"
Set-TextValue "D6" "1.597"
Set-TextValue "E6" "2.97%"
Set-TextValue "D7" "0.9207"
Set-TextValue "E7" "-0.56%"
Set-TextValue "E8" "1.70%"
Set-TextValue "D9" "0.1183"
Set-TextValue "E9" "2.16%"
Set-TextValue "D10" "0.1829"
Set-TextValue "E10" "3.81%"
Set-TextValue "D11" "0.08939"
Set-TextValue "E11" "2.24%"
Set-TextValue "D12" "0.04135"
Set-TextValue "E12" "-5.01%"
Set-TextValue "E13" "-0.18%"
Set-TextValue "D14" "0.001286"
Set-TextValue "E14" "1.61%"
Set-TextValue "D15" "0.005943"
Set-TextValue "E15" "-0.22%"
Set-TextValue "D16" "3.342"
Set-TextValue "E16" "0.02%"
Set-TextValue "E17" "1.95%"
Set-TextValue "E18" "1.43%"
Set-TextValue "D19" "8.309"
Set-TextValue "E19" "6.18%"
Set-TextValue "D20" "0.1352"
Set-TextValue "E20" "-2.77%"
Set-TextValue "D21" "0.3105"
Set-TextValue "E21" "11.92%"
Set-TextValue "D22" "0.04099"
Set-TextValue "E22" "4.54%"
Set-TextValue "D24" "0.003889"
Set-TextValue "E24" "2.46%"
Set-TextValue "E25" "8.30%"
Set-TextValue "D38" "0.02395"
Set-TextValue "E38" "3.88%"
Set-TextValue "D39" "0.05223"
Set-TextValue "E39" "2.99%"
Set-TextValue "D40" "0.006969"
Set-TextValue "E40" "21.24%"
Set-TextValue "D41" "0.007773"
Set-TextValue "E41" "-1.06%"
Set-TextValue "D42" "0.1325"
Set-TextValue "E42" "3.03%"
Set-TextValue "D43" "0.007410"
Set-TextValue "E43" "0.16%"
Set-TextValue "D44" "0.007121"
Set-TextValue "E44" "-2.04%"
Set-TextValue "D45" "0.2981"
Set-TextValue "E45" "-6.29%"
Set-TextValue "D46" "0.00006585"
Set-TextValue "E46" "6.36%"
Set-TextValue "E47" "-0.04%"
Set-TextValue "D48" "0.04538"
Set-TextValue "E48" "-5.46%"
Set-TextValue "D49" "0.004206"
Set-TextValue "E49" "0.11%"
Set-TextValue "E50" "-0.04%"
Set-TextValue "D51" "0.0002002"
Set-TextValue "E51" "-0.04%"
